# Update dSF (column F) values on Sheet1 per repulled/pushed data and mean calculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    11 = 2
    16 = -3
    17 = -3
    19 = 5
    24 = 1
    25 = -3
    26 = 2
    28 = -3
    31 = -1
    35 = 4
    36 = 1
    37 = -4
    38 = -4
    42 = 1
    43 = 1
    50 = 2
    53 = -4
    55 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
